$wb = $excel.ActiveWorkbook

# --- Sheet "保險" (insurance), the 7th worksheet ---
$ws7 = $wb.Worksheets.Item(7)

# Extend formatting for the new header cells (E1:K1) to match the existing
# header style used by B1:D1, then fill in the header labels. The schema
# mirrors the other sheets: company, name, owner, property_category,
# category, date, legislator_name, legislator_id, source_file, index.
$ws7.Range("D1").Copy()
$ws7.Range("E1:K1").PasteSpecial(-4122)  # xlPasteFormats
$ws7.Range("B1").Value = "company"
$ws7.Range("C1").Value = "name"
$ws7.Range("D1").Value = "owner"
$ws7.Range("E1").Value = "property_category"
$ws7.Range("F1").Value = "category"
$ws7.Range("G1").Value = "date"
$ws7.Range("H1").Value = "legislator_name"
$ws7.Range("I1").Value = "legislator_id"
$ws7.Range("J1").Value = "source_file"
$ws7.Range("K1").Value = "index"

# Extend formatting for the new data cells (E2:K2) to match the existing
# data-row style used by B2:D2, then fill in the record's values.
$ws7.Range("D2").Copy()
$ws7.Range("E2:K2").PasteSpecial(-4122)  # xlPasteFormats
$ws7.Range("B2").Value = "國寶人壽"
$ws7.Range("C2").Value = "得意年年"
$ws7.Range("D2").Value = "何欣純"
$ws7.Range("E2").Value = "insurance"
$ws7.Range("F2").Value = "normal"
# "date"-like text must stay text, not become an Excel date serial: force
# the Text format before writing, then restore the original (border-less)
# format by pasting it back in from a sibling data cell.
$ws7.Range("G2").NumberFormat = "@"
$ws7.Range("G2").Value = "2012-04-30"
$ws7.Range("D2").Copy()
$ws7.Range("G2").PasteSpecial(-4122)  # xlPasteFormats
$ws7.Range("H2").Value = "何欣純"
$ws7.Range("I2").Value = 1733
$ws7.Range("J2").Value = "tmp2e891"
$ws7.Range("K2").Value = 80

# --- Sheet "債務" (debt), the 8th worksheet ---
$ws8 = $wb.Worksheets.Item(8)

# Extend formatting for the new header cells (H1:N1), then fill labels.
$ws8.Range("G1").Copy()
$ws8.Range("H1:N1").PasteSpecial(-4122)  # xlPasteFormats
$ws8.Range("B1").Value = "species"
$ws8.Range("C1").Value = "debtor"
$ws8.Range("D1").Value = "owner"
$ws8.Range("E1").Value = "total"
$ws8.Range("F1").Value = "register_date"
$ws8.Range("G1").Value = "register_reason"
$ws8.Range("H1").Value = "property_category"
$ws8.Range("I1").Value = "category"
$ws8.Range("J1").Value = "date"
$ws8.Range("K1").Value = "legislator_name"
$ws8.Range("L1").Value = "legislator_id"
$ws8.Range("M1").Value = "source_file"
$ws8.Range("N1").Value = "index"

# Extend formatting for the new data cells (H2:N2), then fill record values.
$ws8.Range("G2").Copy()
$ws8.Range("H2:N2").PasteSpecial(-4122)  # xlPasteFormats
$ws8.Range("B2").Value = "房屋貸款"
$ws8.Range("C2").Value = "謝俊雄"
$ws8.Range("D2").Value = "霧峰鄕農會臺中市霧峰區四德路"
$ws8.Range("E2").Value = 6081426
$ws8.Range("F2").Value = "93年09月14日"
$ws8.Range("G2").Value = "房貸"
$ws8.Range("H2").Value = "debt"
$ws8.Range("I2").Value = "normal"
# Same text-not-date treatment for the "date" column here too.
$ws8.Range("J2").NumberFormat = "@"
$ws8.Range("J2").Value = "2012-04-30"
$ws8.Range("G2").Copy()
$ws8.Range("J2").PasteSpecial(-4122)  # xlPasteFormats
$ws8.Range("K2").Value = "何欣純"
$ws8.Range("L2").Value = 1733
$ws8.Range("M2").Value = "tmp2e891"
$ws8.Range("N2").Value = 90

# --- Sheet "具有相當價值之財產" (property of considerable value) ---
# The shared-string "otherbonds" is renamed to "antique" in this edit; the
# only cell that referenced it was F2 on this sheet.
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("F2").Value = "antique"
